# Apply the edit described in the commit: changed pi to π
# In this workbook that change is represented as the label in cell B17
# of Sheet1 being renamed from ":special" to ":astral".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the label text
$ws.Range("B17").Value = ":astral"

# Move the active selection to B24, matching the author's final cursor position
$ws.Range("B24").Select()
